$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the typo "ned" -> "end" in the paragraph that reads:
#    "...still items in the list, add the item indexed to the new
#     string, from the ned //to the beginning."
#    The final text should read "...from the end //to the beginning."
#    and the word break left behind (marked with the _GoBack bookmark,
#    as Word does for the most recent edit point) sits between the
#    "en" and the "d" of "end".
# ------------------------------------------------------------------

# Locate the paragraph so we only ever search/replace inside it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*still items in the list*") {
        $target = $cand
        break
    }
}

$pStart = $target.Range.Start
$pText = $target.Range.Text

# Temporary bookmarks placed exactly on the two run boundaries that must
# be preserved ("...new string" | ", from the ..." ) keep those runs
# from being re-coalesced by the edit below.
$bnd1 = $pStart + $pText.IndexOf("still items in the list")
$d.Bookmarks.Add("zzTmpBoundaryA", $d.Range($bnd1, $bnd1))

$bnd2 = $pStart + $pText.IndexOf(", from the")
$d.Bookmarks.Add("zzTmpBoundaryB", $d.Range($bnd2, $bnd2))

# Replace ", from the ned " with ", from the end " in one go so the
# stray proofErr (gramStart/gramEnd) wrapping "ned" is cleared away.
$full = $d.Content.Text
$oldChunk = ", from the ned "
$newChunk = ", from the end "
$idx = $full.IndexOf($oldChunk)
$r = $d.Range($idx, $idx + $oldChunk.Length)
$r.Text = $newChunk

# Remove the helper bookmarks now that the edit is done.
$d.Bookmarks.Item("zzTmpBoundaryA").Delete()
$d.Bookmarks.Item("zzTmpBoundaryB").Delete()

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: delete it from its old spot (between
#    "4" and "n+4" in the "Run time:  4n+4" line) and drop it between
#    "en" and "d" of the word we just fixed above.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$full2 = $d.Content.Text
$endIdx = $full2.IndexOf("end //to the beginning")
$splitPos = $endIdx + 2
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))
